$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.181.07"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.23%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.600.67"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.13%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.000"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "302.92"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.52%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3780"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.65"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3610"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.261"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08112"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.54"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.576"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.377"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001246"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.85%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.602.76"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.61"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06860"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.00"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.526"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("E23").Value = "  -0.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.189.73"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.394"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.987"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +9.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.15"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.12"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.230"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.58"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.18%  "

$ws.Range("E31").Value = "  -0.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.817"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.778.98"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9768"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.86%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07537"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.31"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.60%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02719"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.126"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2500"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.80%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08785"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7087"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.22%  "

$ws.Range("E42").Value = "  -2.01%  "

$ws.Range("E43").Value = "  -2.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.45"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6532"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.302"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.016"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.99%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.15"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07958"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.202"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.230"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.76%  "
